# Adds three new bibliography entries (Pedregosa/scikit-learn, Van Rossum/Python,
# Di Tommaso/Nextflow) to the "Sheet1" bibliography list, continuing the existing
# list style (cell style index "2" == Arial 8pt FF222222) and rich-text citation
# formatting (italic journal/container titles) used throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0

# Arial 8pt #FF222222 is the font baked into cell style "2" (used by every other
# bibliography row), and is also the explicit run font used on every formatted
# run after the first one throughout this workbook's existing citations.
$runFontName = "Arial"
$runFontSize = 8
$runFontColor = 0x222222  # BGR-packed RGB integer Excel COM expects (R=G=B=0x22 here)

# --- Row 93: Pedregosa et al., scikit-learn --------------------------------
$ws.Range("A92").Copy() | Out-Null
$ws.Range("A93").PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> copies style "2"

$text93 = "Pedregosa, F., Varoquaux, G., Gramfort, A., Michel, V., Thirion, B., Grisel, O., ... & Vanderplas, J. (2011). Scikit-learn: Machine learning in Python.$nbsp" + `
          "Journal of Machine Learning Research" + `
          ",$nbsp" + `
          "12" + `
          "(Oct), 2825-2830."
$ws.Range("A93").Value = $text93

$c93 = $ws.Range("A93")

$r1 = $c93.Characters(153, 36)   # "Journal of Machine Learning Research"
$r1.Font.Italic = $true
$r1.Font.Name = $runFontName
$r1.Font.Size = $runFontSize
$r1.Font.Color = $runFontColor

$r2 = $c93.Characters(189, 2)    # ",<nbsp>"
$r2.Font.Italic = $false
$r2.Font.Name = $runFontName
$r2.Font.Size = $runFontSize
$r2.Font.Color = $runFontColor

$r3 = $c93.Characters(191, 2)    # "12"
$r3.Font.Italic = $true
$r3.Font.Name = $runFontName
$r3.Font.Size = $runFontSize
$r3.Font.Color = $runFontColor

$r4 = $c93.Characters(193, 17)   # "(Oct), 2825-2830."
$r4.Font.Italic = $false
$r4.Font.Name = $runFontName
$r4.Font.Size = $runFontSize
$r4.Font.Color = $runFontColor

# --- Row 94: Van Rossum, Python Programming Language ------------------------
$ws.Range("A92").Copy() | Out-Null
$ws.Range("A94").PasteSpecial(-4122) | Out-Null

$text94 = "Van Rossum, G. (2007, June). Python Programming Language. In$nbsp" + `
          "USENIX Annual Technical Conference" + `
          "$nbsp(Vol. 41, p. 36)."
$ws.Range("A94").Value = $text94

$c94 = $ws.Range("A94")

$r5 = $c94.Characters(62, 34)    # "USENIX Annual Technical Conference"
$r5.Font.Italic = $true
$r5.Font.Name = $runFontName
$r5.Font.Size = $runFontSize
$r5.Font.Color = $runFontColor

$r6 = $c94.Characters(96, 18)    # "<nbsp>(Vol. 41, p. 36)."
$r6.Font.Italic = $false
$r6.Font.Name = $runFontName
$r6.Font.Size = $runFontSize
$r6.Font.Color = $runFontColor

# --- Row 95: Di Tommaso et al., Nextflow (plain text, no rich formatting) ---
$ws.Range("A92").Copy() | Out-Null
$ws.Range("A95").PasteSpecial(-4122) | Out-Null

$ws.Range("A95").Value = "Di Tommaso, P., Chatzou, M., Baraja, P. P., & Notredame, C. (2014). A novel tool for highly scalable computational pipelines."

$excel.CutCopyMode = $false

# --- Move the view roughly where the diff leaves it -------------------------
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 82
$ws.Range("A100").Select() | Out-Null
